$d = $word.ActiveDocument

$find = "Сазвежђе сазвежђе Персеј 2022: 16-25 јануар, 7-16 новембар, 6-15 децембар"
$replace = "Сазвежђе сазвежђе Персеј током 2022. године посматрамо 16-25 јануар, 7-16 новембар, 6-15 децембар"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
